$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CasesTab Neo4j query (cell B2): the Cohort column is dropped
# from the RETURN clause (coalesce(co.cohort_description,'') AS Cohort).
$ws.Range("B2").Value = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n`nMATCH (c)<--(diag:diagnosis)`n MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`n`tWHERE s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in ['T3N0M1', 'T3N1M0'] OPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

# With one less line of text the cell's wrapped row now renders shorter;
# match the resulting row height.
$ws.Rows.Item(2).RowHeight = 290

# Selection moved from B4 to B2.
$ws.Range("B2").Select()
